# Natmi following Dr Hou advice
# Update the LR-pair stats for Snca-Lag3 (rows 2-5) with recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.3950880000000001
$ws.Range("H2").Value = 1.185264
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 6.840255666666667
$ws.Range("N2").Value = 20.520767
$ws.Range("O2").Value = 0.2326221245729483
$ws.Range("P2").Value = 0.2326221245729483
$ws.Range("Q2").Value = 2.702502930832
$ws.Range("R2").Value = 24.322526377488
$ws.Range("S2").Value = 0.2326221245729483
$ws.Range("T2").Value = 0.2326221245729483

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.3950880000000001
$ws.Range("H3").Value = 1.185264
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.415322999999999
$ws.Range("N3").Value = 28.245969
$ws.Range("O3").Value = 0.3201945287620894
$ws.Range("P3").Value = 0.3201945287620895
$ws.Range("Q3").Value = 3.719881133424
$ws.Range("R3").Value = 33.478930200816
$ws.Range("S3").Value = 0.3201945287620894
$ws.Range("T3").Value = 0.3201945287620895

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.3950880000000001
$ws.Range("H4").Value = 1.185264
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.720817
$ws.Range("N4").Value = 23.162451
$ws.Range("O4").Value = 0.2625680883144773
$ws.Range("P4").Value = 0.2625680883144773
$ws.Range("Q4").Value = 3.050402146896
$ws.Range("R4").Value = 27.453619322064
$ws.Range("S4").Value = 0.2625680883144773
$ws.Range("T4").Value = 0.2625680883144773

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.3950880000000001
$ws.Range("H5").Value = 1.185264
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.428613333333334
$ws.Range("N5").Value = 16.28584
$ws.Range("O5").Value = 0.1846152583504849
$ws.Range("P5").Value = 0.1846152583504849
$ws.Range("Q5").Value = 2.14477998464
$ws.Range("R5").Value = 19.30301986176
$ws.Range("S5").Value = 0.1846152583504849
$ws.Range("T5").Value = 0.1846152583504849
